$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The BOM row that used to list "U11, U12" together (row 50) with a
# single part (TPS3703A7120DSERQ1) is being split into two separate
# rows: one for U11 (new part TPS3703A7330DSERQ1) and one for U12
# (keeps the original part TPS3703A7120DSERQ1). This inserts a new
# row (new row 51) and pushes everything below down by one.
# -----------------------------------------------------------------

# Insert a new row right after row 50; it inherits formatting from
# the row above (row 50), matching the target layout.
$ws.Rows("51:51").Insert()

# Row 50 keeps designator U11, but now points to the new part.
$ws.Range("A50").Value = "U11"
$ws.Range("B50").Value = "TPS3703A7330DSERQ1"
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = "https://www.mouser.fr/ProductDetail/Texas-Instruments/TPS3703A7330DSERQ1?qs=%2Fha2pyFadugGt5dAWQFWLSoN%2FsbubeOHv61%2FdlTbBouvNCowEjdZ9Q%3D%3D"

# New row 51 holds designator U12 with the original part.
$ws.Range("A51").Value = "U12"
$ws.Range("B51").Value = "TPS3703A7120DSERQ1"
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = "https://www.mouser.fr/ProductDetail/Texas-Instruments/TPS3703A7120DSERQ1?qs=%2Fha2pyFadui8P4qqv7wtJq7XkkCwOVhRNdsWZTFkD%252B34OV3EnMk3KZyek4uztLX0"

# -----------------------------------------------------------------
# Rebuild all hyperlinks on the sheet. Inserting a row does not shift
# the existing Hyperlinks collection in this engine, so the safest
# route is to wipe them all and re-add them in the correct, final
# row positions.
# -----------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range('D10'), 'https://fr.farnell.com/tdk/c1608x7s2a473k080ab/cond-0-047-f-100v-10-x7s-0603/dp/2906791?ost=2906791')
$ws.Hyperlinks.Add($ws.Range('D12'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/KRM55WR71H336MH01L?qs=QzBtWTOodeWVwA4bjpST6w%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM033R71E103KE14D?qs=%2Fha2pyFaduiG1IarW8zyNYGyruWgGg3iGVRKpsjt43bog7Vdn557JJK3OhmiUPBA')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM033C71A104KE14D?qs=%2Fha2pyFaduiIqi3iE1rCfs0%2F1GruPdkSQPaytBVK3Ye62TVOf67W82YiR%252BmyGnM0')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://www.mouser.fr/ProductDetail/TDK/C3225X7R1N106K250AC?qs=P1JMDcb91o5%2Fi%252BA2Vs0UoQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D5'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM188R72A104KA35D?qs=P%252BBA3F6RM%2F5aaC9LFKRzYw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D6'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM033C81A105ME05D?qs=%2Fha2pyFaduiC0uiaBjWXb2K8VBBqJ%252BZm09tmNRRx24dUCYe0%2F77QZfTbjVNdfe%2Fp')
$ws.Hyperlinks.Add($ws.Range('D7'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM188R61A226ME15D?qs=eeBpzGFlv%252B8DV%2FrilzyhAw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D8'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM188R61A106ME69D?qs=o98hbGm2QKGN31NsNZp5HQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D9'), 'https://www.mouser.fr/ProductDetail/TDK/C1005X5R1V105K050BE?qs=%2Fha2pyFaduhU5KeTuzJ6cpvEC5KCh3tjzdHX67H3uw%252B2LA1zZVn4jrUs3cYKP86b')
$ws.Hyperlinks.Add($ws.Range('D11'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/GRM033R71H102KA12D?qs=2W5sgKM%2F373OnvLttcxBWw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D13'), 'https://fr.farnell.com/broadcom-limited/hsmg-c280/led-vert-15mcd-572nm-0402/dp/2494322?ost=2494322')
$ws.Hyperlinks.Add($ws.Range('D14'), 'https://fr.farnell.com/broadcom-limited/hsmc-c280/led-rouge-90mcd-626nm-0402/dp/2494321?ost=2494321')
$ws.Hyperlinks.Add($ws.Range('D15'), 'https://www.tme.eu/fr/details/ws2812b-b/diodes-led-smd-couleur/worldsemi/ws2812b-black/')
$ws.Hyperlinks.Add($ws.Range('D17'), 'https://www.mouser.fr/ProductDetail/Murata-Electronics/BLM18SG221TN1D?qs=%2Fha2pyFadug5NtL6n1y3hLpQXfoAMDm5SdpiQ1riGOvcL36yZTLdNw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D16'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/TPD1E04U04DPYR?qs=%2Fha2pyFaduhXCCun1enApVgGMl1srv71hbbstzc3NtDNeFqEPLjqvZTVaXB%252B1BBY')
$ws.Hyperlinks.Add($ws.Range('D18'), 'https://www.mouser.fr/ProductDetail/Hirose-Connector/DF13-11P-125DSA?qs=Ux3WWAnHpjBlXPmXUKf1pg%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D19'), 'https://www.mouser.fr/ProductDetail/Hirose-Connector/DF13-5P-125DSA?qs=Ux3WWAnHpjDJWf8XCjmDFw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D20'), 'https://www.mouser.fr/ProductDetail/Hirose-Connector/DF13-3P-125DSA?qs=Ux3WWAnHpjDQ1kTMJnVTAQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D21'), 'https://www.mouser.fr/ProductDetail/Hirose-Connector/DF13-2P-125DSA?qs=Ux3WWAnHpjA3bAkWKtiTkA%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D22'), 'https://www.mouser.fr/ProductDetail/Amphenol-FCI/20021311-00008T4LF?qs=%2Fha2pyFaduhRt0ldyeOmZUMWVTlP60Btyx%2F7BwI%252BAmDSTHkMifvMxe8P09SyfmPt')
$ws.Hyperlinks.Add($ws.Range('D23'), 'https://www.mouser.fr/ProductDetail/Hirose-Connector/DF13-7P-125DSA?qs=Ux3WWAnHpjDNnwT%252BBKTPkw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D24'), 'https://www.tme.eu/fr/details/mr30pw-fb/connecteurs-dc/amass/')
$ws.Hyperlinks.Add($ws.Range('D25'), 'https://www.tme.eu/fr/details/xt30pw-m/connecteurs-dc/amass/')
$ws.Hyperlinks.Add($ws.Range('D26'), 'https://www.mouser.fr/ProductDetail/Samtec/ADF6-10-035-L-4-2-A-TR?qs=%2Fha2pyFaduj%2F1bGGyxbkubC23W7y4EFX4SCR%2F%252B6eh6LPxktvNuI7OzoMle1t4Dp5')
$ws.Hyperlinks.Add($ws.Range('D27'), 'https://www.mouser.fr/ProductDetail/Coilcraft/LPS6235-223MRC?qs=%2Fha2pyFaduh%252BDrgCJSwu8Dw5nTl6KPd2aJHbjQ7CO3gsK9BVvLJWFFhaeY2409J6')
$ws.Hyperlinks.Add($ws.Range('D28'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/CSD88599Q5DC?qs=5aG0NVq1C4z5gVbcQ76feA%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D29'), 'https://www.mouser.fr/ProductDetail/Nexperia/NX7002BKMYL?qs=%2Fha2pyFaduhoqz%2FLZn5LcHsDOjOn5prphHYCU0gExPHTshUDH3RDsQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D30'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF1200C?qs=%2Fha2pyFadujWqMGonbcVXWXDNnBbNvNsbAVwT70CvkaubqXTVjfaxQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D31'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF1003C?qs=%2Fha2pyFadujWqMGonbcVXUyWU9hohrWJgVqtZdVl5k%2FvHnuychE5ng%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D32'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF2212C?qs=%2Fha2pyFadujWqMGonbcVXZTd6zTO7cT5FlrWhgbZRCq0%2FjoFBD6z6A%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D33'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF33R0C?qs=%2Fha2pyFadujWqMGonbcVXQUm8VVN5JAr6%252BqK2V5ChuscKfrTrS8qrw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D34'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF56R0C?qs=%2Fha2pyFadujWqMGonbcVXXz6Wadq5KM%2FgA6yFLBK%2FJx1N74UoX4rAQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D35'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF1001C?qs=%2Fha2pyFadujWqMGonbcVXSSH66T3l46T3z23R1pKr2zx36TYTYFg%252Bw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D36'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF1000C?qs=%2Fha2pyFadujWqMGonbcVXfzkRhIhrzCTGf59nzrXMdmP0MGzupBPvg%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D37'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF1002C?qs=%2Fha2pyFadujWqMGonbcVXdOjv8PYAhpdlCZ7N6tv00WxAhYaPU3s6g%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D38'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-2RKF1003X?qs=%2Fha2pyFadujWqMGonbcVXa7JUN8iO44iYRctk1GHihjbiPKPpnD1Eg%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D39'), 'https://www.mouser.fr/ProductDetail/Ohmite/FC4L64R005FER?qs=%2Fha2pyFaduhbX7CsPg9tX0CyS2vmh4VF4JaBdYQojvRSt3tE7ez6cQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D40'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-2RKF5101X?qs=%2Fha2pyFadujWqMGonbcVXYbnJN%2FW1vPkd9Ro2q4XkzoQcicelgdF1w%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D41'), 'https://www.mouser.fr/ProductDetail/Panasonic/ERJ-1GNF4701C?qs=%2Fha2pyFadujWqMGonbcVXbI2DeneIvW9q5JbJdNDaELo4QMncwvhrQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D42'), 'https://www.mouser.fr/ProductDetail/Omron-Electronics/B3U-1000P?qs=AO7BQMcsEu4ip80xyf2FwA%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D43'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/F28388DZWTS?qs=mAH9sUMRCtv%2F6yspAHup2w%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D44'), 'https://www.mouser.fr/ProductDetail/Maxim-Integrated/MAX17502FATB%2bT?qs=%2Fha2pyFadughqAITL5wNFuNoLQrtGJSX3XmITHxKAqjaHXU1qiEl83mAFmMjl6DZ')
$ws.Hyperlinks.Add($ws.Range('D45'), 'https://www.mouser.fr/ProductDetail/Monolithic-Power-Systems-MPS/MPM3804GG-Z?qs=%2Fha2pyFaduhvS2goSYr2M%2FKcg8OWFam93oAF6h5cETljIY3Iv0Vk6g%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D46'), 'https://www.mouser.fr/ProductDetail/Analog-Devices/LT6654AIS6-3TRMPBF?qs=%2Fha2pyFaduhR5nrRyokg0I%2F9GFIkREkb5p7iutAs1DI4fa4ie2YHPaJ8nskXE%2FvV')
$ws.Hyperlinks.Add($ws.Range('D47'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/LSF0204RUTR?qs=%2Fha2pyFadujf%2Fkl32k2TO%2Flb0stQCl8MDqxTYLTZ6T7GqVIp7fKgKQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D48'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/TCAN332GDCNR?qs=%2Fha2pyFadugypkkiYbgeDKOKpUyK5GklIqfKQXJEFO%252BcYfmRbCtCBQ%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D49'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/DSLVDS1048PWR?qs=%2Fha2pyFadui%2FTXskH3hgvdyHYEtJOlOyaJiQL0JEY3YjjG58Bx8NMDxqJeNgmMQY')
$ws.Hyperlinks.Add($ws.Range('D50'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/TPS3703A7330DSERQ1?qs=%2Fha2pyFadugGt5dAWQFWLSoN%2FsbubeOHv61%2FdlTbBouvNCowEjdZ9Q%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D52'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/TPD6E05U06RVZR?qs=%2Fha2pyFaduhaHCMZV6EvVM4X8OdiEwU%252BZhpJxs1AFUrqba2DJH%2FLth2ECKd%252B8vcu')
$ws.Hyperlinks.Add($ws.Range('D53'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/DSLVDS1047PWR?qs=%2Fha2pyFadui%2FTXskH3hgvTDdFAWoPqNzm%2FfR6ePjSevLll9umBcJiva2yCWVx2Lh')
$ws.Hyperlinks.Add($ws.Range('D54'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/DRV8353SRTAT?qs=%2Fha2pyFaduilsGqYdRS7eQddawPYear%2FHTAryIjK2dn%2FkhHr7CfvQw%3D%3D')
$ws.Hyperlinks.Add($ws.Range('D55'), 'https://www.mouser.fr/ProductDetail/ABRACON/ASTX-H12-25000MHZ-T?qs=%2Fha2pyFadug5wyn1HeCucjLZIt0DMOdv3ho1L8YNpV%252BjlXs4FJNH82rKisUL7d9u')
$ws.Hyperlinks.Add($ws.Range('D51'), 'https://www.mouser.fr/ProductDetail/Texas-Instruments/TPS3703A7120DSERQ1?qs=%2Fha2pyFadui8P4qqv7wtJq7XkkCwOVhRNdsWZTFkD%252B34OV3EnMk3KZyek4uztLX0')

# -----------------------------------------------------------------
# Update the view state to match where the author left off editing.
# -----------------------------------------------------------------
$ws.Range("B50").Select()
